$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '37.925.22'
Set-TextValue 'E2' '  -0.79%  '
Set-TextValue 'D3' '2.047.24'
Set-TextValue 'E3' '  -0.47%  '
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '228.43'
Set-TextValue 'E5' '  +0.08%  '
Set-TextValue 'E6' '  -1.18%  '
Set-TextValue 'D7' '60.99'
Set-TextValue 'E7' '  +0.65%  '
Set-TextValue 'E8' '  -0.06%  '
Set-TextValue 'D9' '0.377'
Set-TextValue 'E9' '  -2.27%  '
Set-TextValue 'D10' '0.0822'
Set-TextValue 'E10' '  -0.49%  '
Set-TextValue 'E11' '  +0.43%  '
Set-TextValue 'B12' 'Chainlink'
Set-TextValue 'C12' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D12' '14.67'
Set-TextValue 'E12' '  -0.82%  '
Set-TextValue 'B13' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C13' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D13' '2.344.50'
Set-TextValue 'E13' '  -0.79%  '
Set-TextValue 'D14' '21.07'
Set-TextValue 'E14' '  +0.20%  '
Set-TextValue 'D15' '0.778'
Set-TextValue 'E15' '  +2.34%  '
Set-TextValue 'D16' '5.22'
Set-TextValue 'E16' '  -1.49%  '
Set-TextValue 'D17' '2.060.93'
Set-TextValue 'E17' '  +0.32%  '
Set-TextValue 'D18' '37.837.09'
Set-TextValue 'E18' '  -0.82%  '
Set-TextValue 'D19' '69.60'
Set-TextValue 'E19' '  -0.22%  '
Set-TextValue 'D20' '5.92'
Set-TextValue 'E20' '  -4.76%  '
Set-TextValue 'D21' '0.0₃0824'
Set-TextValue 'E21' '  -1.03%  '
Set-TextValue 'D22' '224.18'
Set-TextValue 'E22' '  -0.35%  '
Set-TextValue 'E23' '  +0.09%  '
Set-TextValue 'D24' '2.44'
Set-TextValue 'E24' '  -0.05%  '
Set-TextValue 'D25' '2.27'
Set-TextValue 'E25' '  +2.47%  '
Set-TextValue 'D26' '168.17'
Set-TextValue 'E26' '  +0.98%  '
Set-TextValue 'D27' '9.34'
Set-TextValue 'E27' '  +1.08%  '
Set-TextValue 'E28' '  -2.32%  '
Set-TextValue 'D29' '18.83'
Set-TextValue 'E29' '  -0.71%  '
Set-TextValue 'E30' '  -1.92%  '
Set-TextValue 'E31' '  -0.18%  '
Set-TextValue 'E32' '  +8.43%  '
Set-TextValue 'D33' '4.39'
Set-TextValue 'E33' '  -1.85%  '
Set-TextValue 'D34' '4.51'
Set-TextValue 'E34' '  -0.57%  '
Set-TextValue 'D35' '0.0600'
Set-TextValue 'E35' '  -0.70%  '
Set-TextValue 'D36' '6.61'
Set-TextValue 'E36' '  +4.54%  '
Set-TextValue 'D37' '2.34'
Set-TextValue 'E37' '  +2.72%  '
Set-TextValue 'D38' '3.47'
Set-TextValue 'E38' '  +6.61%  '
Set-TextValue 'E39' '  -0.17%  '
Set-TextValue 'D40' '18.15'
Set-TextValue 'E40' '  +8.02%  '
Set-TextValue 'D41' '1.540.50'
Set-TextValue 'E41' '  +0.96%  '
Set-TextValue 'E42' '  -0.27%  '
Set-TextValue 'D43' '96.43'
Set-TextValue 'E43' '  -1.29%  '
Set-TextValue 'E44' '  -1.02%  '
Set-TextValue 'D45' '0.0912'
Set-TextValue 'E45' '  -1.65%  '
Set-TextValue 'D46' '4.18'
Set-TextValue 'E46' '  +4.12%  '
Set-TextValue 'E47' '  -1.46%  '
Set-TextValue 'E48' '  -0.28%  '
Set-TextValue 'E49' '  -1.09%  '
Set-TextValue 'D50' '7.06'
Set-TextValue 'E50' '  +0.03%  '
Set-TextValue 'D51' '2.236.02'
